# Opsplitsen per subset ingebouwd.
# TODO: Significantie tegenover een totaal van een subset toevoegen.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "onderdelen" sheet: split the "GGD NOG" / 2022 configuration into two
# subsets (Gemeentecode, subregio) and tag the "NOG 2019" row with its
# own subset too.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("onderdelen")

# Make room: insert a new row above the old row 3 ("NOG 2019").
$ws.Rows.Item(3).Insert()

# Row 2 (GGD NOG / 2022 / TRUE / intern / NL totaal) now belongs to the
# "Gemeentecode" subset.
$ws.Range("B2").Value = "Gemeentecode"

# New row 3: a second GGD NOG / 2022 configuration for the "subregio" subset.
$ws.Range("A3").Value = "GGD NOG"
$ws.Range("B3").Value = "subregio"
$ws.Range("C3").Value = 2022
$ws.Range("D3").Value = $false

# Row 4 (previously row 3: "NOG 2019" / 2019) belongs to the "Gemeentecode"
# subset as well.
$ws.Range("B4").Value = "Gemeentecode"

# Widen the new "subset" column (B) to fit its contents.
$ws.Columns.Item(2).ColumnWidth = 20.166666666666668

# ------------------------------------------------------------------
# "datasets" sheet: widen columns D/E and move the remembered selection,
# without changing which sheet/tab is active.
# ------------------------------------------------------------------
$wsDatasets = $wb.Worksheets.Item("datasets")
$wsDatasets.Activate() | Out-Null
$wsDatasets.Columns.Item(4).ColumnWidth = 26.666666666666668
$wsDatasets.Columns.Item(5).ColumnWidth = 22.5
$wsDatasets.Range("D2").Select() | Out-Null

# Restore "onderdelen" as the active/visible sheet.
$ws.Activate() | Out-Null
